# Update "想去人数" (interest count) figures in the F column across the
# "展览", "演出" and "全部类型" sheets to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 1712
$ws.Range("F3").Value  = 9199
$ws.Range("F4").Value  = 122
$ws.Range("F6").Value  = 723
$ws.Range("F7").Value  = 1403
$ws.Range("F8").Value  = 213
$ws.Range("F9").Value  = 73
$ws.Range("F10").Value = 110
$ws.Range("F11").Value = 6008
$ws.Range("F12").Value = 62
$ws.Range("F13").Value = 394
$ws.Range("F15").Value = 4777
$ws.Range("F16").Value = 21
$ws.Range("F17").Value = 171
$ws.Range("F18").Value = 1156
$ws.Range("F19").Value = 43
$ws.Range("F21").Value = 44
$ws.Range("F22").Value = 3
$ws.Range("F23").Value = 273
$ws.Range("F25").Value = 3226
$ws.Range("F26").Value = 134

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 57

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 1712
$ws.Range("F3").Value  = 9199
$ws.Range("F4").Value  = 122
$ws.Range("F5").Value  = 57
$ws.Range("F7").Value  = 723
$ws.Range("F8").Value  = 1403
$ws.Range("F9").Value  = 213
$ws.Range("F10").Value = 73
$ws.Range("F11").Value = 110
$ws.Range("F12").Value = 6008
$ws.Range("F13").Value = 62
$ws.Range("F14").Value = 394
$ws.Range("F16").Value = 4777
$ws.Range("F17").Value = 21
$ws.Range("F18").Value = 171
$ws.Range("F19").Value = 1156
$ws.Range("F20").Value = 43
$ws.Range("F22").Value = 44
$ws.Range("F23").Value = 3
$ws.Range("F24").Value = 273
$ws.Range("F26").Value = 3226
$ws.Range("F28").Value = 134
